$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Illustration Calculator")
$ws2 = $wb.Worksheets.Item("Plan Calculations")

# Change the selected plan from "Plan #4" to "Plan #2"
$ws1.Range("B6").Value = "Plan #2"
# Change current age of client from 40 to 68
$ws1.Range("B9").Value = 68

# Fix the Plan #2 payout formula in column G of "Plan Calculations":
# once the client's starting age (C3, i.e. B9 on Illustration Calculator)
# is already >= 65, compute the payout off of the payment-year index (D)
# rather than off of the (now irrelevant) age-minus-65 logic.
$g3Formula = "=IF(R3C3>=65,IF(RC[-4]>94,0,25000+((RC[-3]-1)*1000)),IF(OR(RC[-4]<65,RC[-4]>94),0,25000+((RC[-4]-65)*1000)))"
$ws2.Range("G3").FormulaR1C1 = $g3Formula
$ws2.Range("G4:G67").FormulaR1C1 = $g3Formula
$ws2.Range("G68:G75").FormulaR1C1 = $g3Formula

# Restore cell selections to match the saved state of the workbook
[void]$ws2.Activate()
[void]$ws2.Range("G4").Select()
[void]$ws1.Activate()
[void]$ws1.Range("B10").Select()
